$d = $word.ActiveDocument

$old = '<w:style w:type="paragraph" w:styleId="BlockText"><w:name w:val="Block Text"/><w:basedOn w:val="BodyText"/><w:next w:val="BodyText"/><w:uiPriority w:val="9"/><w:unhideWhenUsed/><w:qFormat/><w:pPr><w:spacing w:before="100" w:after="100"/><w:ind w:firstLine="0"/></w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:style>'
$new = '<w:style w:type="paragraph" w:styleId="BlockText"><w:name w:val="Block Text"/><w:basedOn w:val="BodyText"/><w:next w:val="BodyText"/><w:uiPriority w:val="9"/><w:unhideWhenUsed/><w:qFormat/><w:pPr><w:spacing w:before="100" w:after="100"/><w:ind w:firstLine="0" w:left="480" w:right="480"/></w:pPr></w:style>'

$xml = $d.WordOpenXML
if ($xml.Contains($old)) {
    $xml = $xml.Replace($old, $new)
    $d.WordOpenXML = $xml
    Write-Output "replaced"
} else {
    Write-Output "PATTERN NOT FOUND"
}
